# Actualiza los horarios de la Linea 141 (commit: "Horarios actualizados Linea 141 - 562")
# Reescribe las 3 hojas (LP1912, LP1912-215, 6203-6173) con los datos nuevos scrapeados
# a las 06:34:35, agregando las filas nuevas y manteniendo el resto ordenado por Hora_Llegada.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Hoja 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 06:34:35"
$ws1.Range("A3").Value = "Total filas: 30"

$rows1 = @(
    @("05:57:04","06:09","10_OLMOS",12,"LP1912"),
    @("05:57:04","06:16","215A_EL PATO",19,"LP1912"),
    @("05:57:04","06:30","23_HERNANDEZ",33,"LP1912"),
    @("05:57:04","06:34","11_ETCHEVERRY",37,"LP1912"),
    @("06:34:35","06:36","11_ETCHEVERRY",2,"LP1912"),
    @("05:57:04","06:39","17X38_ROMERO",42,"LP1912"),
    @("05:57:04","06:41","16_SANTA ANA",44,"LP1912"),
    @("06:16:41","06:56","215A_EL PATO",40,"LP1912"),
    @("05:57:04","06:57","215A_EL PATO",60,"LP1912"),
    @("05:57:04","06:59","225_GOMEZ",62,"LP1912"),
    @("06:16:41","07:15","215C_EL PATO",59,"LP1912"),
    @("05:57:04","07:16","215C_EL PATO",79,"LP1912"),
    @("05:57:04","07:19","14_ABASTO",82,"LP1912"),
    @("06:16:41","07:20","16_SANTA ANA",64,"LP1912"),
    @("06:16:41","07:21","23_HERNANDEZ",65,"LP1912"),
    @("05:57:04","07:21","16_SANTA ANA",84,"LP1912"),
    @("05:57:04","07:22","23_HERNANDEZ",85,"LP1912"),
    @("05:57:04","07:29","17X38_ROMERO",92,"LP1912"),
    @("05:57:04","07:35","10_OLMOS",98,"LP1912"),
    @("06:16:41","07:36","27_EL RETIRO",80,"LP1912"),
    @("05:57:04","07:37","27_EL RETIRO",100,"LP1912"),
    @("06:34:35","07:43","215A_EL PATO",69,"LP1912"),
    @("05:57:04","07:55","14_ABASTO",118,"LP1912"),
    @("06:16:41","08:00","17_ROMERO",104,"LP1912"),
    @("06:16:41","08:01","16_SANTA ANA",105,"LP1912"),
    @("06:34:35","08:06","23_HERNANDEZ",92,"LP1912"),
    @("06:16:41","08:11","10_OLMOS",115,"LP1912"),
    @("06:16:41","08:13","15X38_ABASTO",117,"LP1912"),
    @("06:34:35","08:29","11_ETCHEVERRY",115,"LP1912"),
    @("06:34:35","08:29","15_ABASTO",115,"LP1912")
)

$r = 6
foreach ($row in $rows1) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Hoja 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 06:34:35"
$ws2.Range("A3").Value = "Total filas: 6"

$rows2 = @(
    @("05:57:04","06:16","215A_EL PATO",19,"LP1912"),
    @("06:16:41","06:56","215A_EL PATO",40,"LP1912"),
    @("05:57:04","06:57","215A_EL PATO",60,"LP1912"),
    @("06:16:41","07:15","215C_EL PATO",59,"LP1912"),
    @("05:57:04","07:16","215C_EL PATO",79,"LP1912"),
    @("06:34:35","07:43","215A_EL PATO",69,"LP1912")
)

$r = 6
foreach ($row in $rows2) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Hoja 3: 6203-6173 (solo cambia el timestamp de actualizacion)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 06:34:35"
